# Refactor the "CI" header row:
#  - keep item/marca/ref/ncm/cor headers (A1:E1)
#  - replace the old bra_dummy + 34..44 + cajas text headers (F1:R1) with a
#    plain numeric size run 20..45 spanning F1:AE1
#  - move/re-label the trailing summary columns to AF1:AI1
#    (cajas -> caixas, total pares, unit price -> preco unit, valor total)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1:E1 - unchanged labels
$ws.Cells.Item(1, 1).Value = "item"
$ws.Cells.Item(1, 2).Value = "marca"
$ws.Cells.Item(1, 3).Value = "ref"
$ws.Cells.Item(1, 4).Value = "ncm"
$ws.Cells.Item(1, 5).Value = "cor"

# F1:AE1 - numeric size run 20..45 (26 columns)
$sizes = 20..45
for ($i = 0; $i -lt $sizes.Length; $i++) {
    $ws.Cells.Item(1, 6 + $i).Value = $sizes[$i]
}

# AF1:AI1 - trailing summary columns (renamed / relocated)
$ws.Cells.Item(1, 32).Value = "caixas"
$ws.Cells.Item(1, 33).Value = "total pares"
$ws.Cells.Item(1, 34).Value = "preco unit"
$ws.Cells.Item(1, 35).Value = "valor total"

# Narrow the new size columns (F:AE) to width 5
$ws.Range("F1:AE1").ColumnWidth = 4.2

# Match the saved selection state
[void]$ws.Range("AM20").Select()
